$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44214
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 1800
$ws.Range("O2").Value = 1800
$ws.Range("P2").Value = 1800
$ws.Range("S2").Value = 1800

# Row 3
$ws.Range("D3").Value2 = 44176
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value2 = 44567
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 2400
$ws.Range("O4").Value = 2400
$ws.Range("P4").Value = 2400
$ws.Range("S4").Value = 2400

# Row 7
$ws.Range("D7").Value2 = 44998
$ws.Range("M7").Value = 20

# Row 8
$ws.Range("D8").Value2 = 44999
$ws.Range("M8").Value = 25
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("S8").Value = 2500

# Row 9
$ws.Range("D9").Value2 = 44175
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 5000
$ws.Range("O9").Value = 5000
$ws.Range("P9").Value = 5000
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 5000

# Row 10
$ws.Range("D10").Value2 = 44323
$ws.Range("N10").Value = 3200
$ws.Range("O10").Value = 3200
$ws.Range("P10").Value = 3200
$ws.Range("S10").Value = 3200

# Row 11
$ws.Range("D11").Value2 = 44215
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = 2800
$ws.Range("O11").Value = 2800
$ws.Range("P11").Value = 2800
$ws.Range("R11").Value = "Región de La Araucanía"
$ws.Range("S11").Value = 2800

# Row 12
$ws.Range("D12").Value2 = 44592
$ws.Range("M12").Value = 5
$ws.Range("N12").Value = 7500
$ws.Range("O12").Value = 7500
$ws.Range("P12").Value = 7500
$ws.Range("R12").Value = "Región de La Araucanía"
$ws.Range("S12").Value = 7500

# Row 13
$ws.Range("D13").Value2 = 44574
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 3000
$ws.Range("O13").Value = 3000
$ws.Range("P13").Value = 3000
$ws.Range("S13").Value = 3000
